$d = $word.ActiveDocument

function Get-ParaPieces($rng) {
    # rng must be a Range positioned (Find) inside the target paragraph.
    $para = $rng.Paragraphs(1)
    $prng = $para.Range.Duplicate
    $pxml = $prng.WordOpenXML
    $ptag = ""
    $pprblock = ""
    if ($pxml -match '(?s)<w:body>(<w:p\b[^>]*>)') { $ptag = $matches[1] }
    if ($pxml -match '(?s)<w:body>.*?(<w:pPr>.*?</w:pPr>)') { $pprblock = $matches[1] }
    # Strip synthesized w14:paraId / w14:textId attributes that are not
    # present in the original document so we do not introduce them.
    $ptag = $ptag -replace ' w14:paraId="[0-9A-Fa-f]*"', ''
    $ptag = $ptag -replace ' w14:textId="[0-9A-Fa-f]*"', ''
    return @{ Range = $prng; Tag = $ptag; PPr = $pprblock }
}

function New-FlatOpcXml($bodyInner) {
    return '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyInner + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# ---------------------------------------------------------------------
# 1) "Setps" -> wrap the existing run with proofErr spellStart/spellEnd
# ---------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("Setps")
$pieces1 = Get-ParaPieces $rng1
$pieces1.Range.MoveEnd(1, -1)
$pieces1.Range.Delete()

$inner1 = '<w:proofErr w:type="spellStart"/><w:r><w:t>Setps</w:t></w:r><w:proofErr w:type="spellEnd"/>'
$full1 = $pieces1.Tag + $pieces1.PPr + $inner1 + '</w:p>'
$pieces1.Range.InsertXML((New-FlatOpcXml $full1))

# ---------------------------------------------------------------------
# 2) "No priviledge " -> split run, wrap "priviledge" with proofErr
#    (keeps remaining runs of the paragraph untouched)
# ---------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("No priviledge ")
$pieces2 = Get-ParaPieces $rng2
$pieces2.Range.MoveEnd(1, -1)
$pieces2.Range.Delete()

$inner2 = '<w:r><w:t xml:space="preserve">No </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>priviledge</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
          '<w:r><w:t>E101 (</w:t></w:r>' + `
          '<w:r><w:t>Not have Authentication</w:t></w:r>' + `
          '<w:r><w:t>)</w:t></w:r>'
$full2 = $pieces2.Tag + $pieces2.PPr + $inner2 + '</w:p>'
$pieces2.Range.InsertXML((New-FlatOpcXml $full2))

# ---------------------------------------------------------------------
# 3) "Please enter correct creditional E112" -> split run, wrap
#    "creditional" with proofErr, and add a new paragraph
#    "User not login before E113" (inheriting the bookmark) after it.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("Please enter correct creditional E112")
$pieces3 = Get-ParaPieces $rng3
$pieces3.Range.Delete()

$inner3 = '<w:r><w:t xml:space="preserve">Please enter correct </w:t></w:r>' + `
          '<w:proofErr w:type="spellStart"/>' + `
          '<w:r><w:t>creditional</w:t></w:r>' + `
          '<w:proofErr w:type="spellEnd"/>' + `
          '<w:r><w:t xml:space="preserve"> E112</w:t></w:r>'
$p1 = $pieces3.Tag + $pieces3.PPr + $inner3 + '</w:p>'

$inner4 = '<w:r><w:t>User not login before E113</w:t></w:r>' + `
          '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
          '<w:bookmarkEnd w:id="0"/>'
$p2 = '<w:p>' + $pieces3.PPr + $inner4 + '</w:p>'

$pieces3.Range.InsertXML((New-FlatOpcXml ($p1 + $p2)))
